$wb = $excel.ActiveWorkbook

# --- Sheet 1: LP1912 ---
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 1).Value = 'Última actualización: 07:31:43'
$ws.Cells.Item(3, 1).Value = 'Total filas: 74'

$rows1 = @(
    @('03:58:57', '04:01', '81_EL PELIGRO', 3, 'LP1912'),
    @('03:58:57', '04:45', '215A_EL PATO', 47, 'LP1912'),
    @('04:28:33', '04:46', '215A_EL PATO', 18, 'LP1912'),
    @('03:58:57', '04:53', '11_ETCHEVERRY', 55, 'LP1912'),
    @('03:58:57', '05:16', '17_ROMERO', 78, 'LP1912'),
    @('03:58:57', '05:21', '23_HERNANDEZ', 83, 'LP1912'),
    @('04:28:33', '05:22', '23_HERNANDEZ', 54, 'LP1912'),
    @('03:58:57', '05:34', '215B_EL PATO', 96, 'LP1912'),
    @('05:31:23', '05:35', '215B_EL PATO', 4, 'LP1912'),
    @('03:58:57', '05:46', '15_ABASTO', 108, 'LP1912'),
    @('03:58:57', '05:53', '10_OLMOS', 115, 'LP1912'),
    @('04:28:33', '05:54', '10_OLMOS', 86, 'LP1912'),
    @('06:01:37', '06:01', '81_EL PELIGRO', 0, 'LP1912'),
    @('04:28:33', '06:04', '16_SANTA ANA', 96, 'LP1912'),
    @('04:28:33', '06:11', '215A_EL PATO', 103, 'LP1912'),
    @('04:28:33', '06:14', '225_HARAS DEL SUR', 106, 'LP1912'),
    @('04:28:33', '06:21', '26_HERNANDEZ', 113, 'LP1912'),
    @('04:28:33', '06:27', '23_HERNANDEZ', 119, 'LP1912'),
    @('04:48:44', '06:29', '86_EST CHICA-ESC AGRARIA', 101, 'LP1912'),
    @('04:48:44', '06:31', '16_SANTA ANA', 103, 'LP1912'),
    @('04:48:44', '06:44', '225_C ROCA-H SUR', 116, 'LP1912'),
    @('04:48:44', '06:46', '215C_EL PATO', 118, 'LP1912'),
    @('05:31:23', '06:59', '14_ABASTO', 88, 'LP1912'),
    @('06:44:15', '07:00', '14_ABASTO', 16, 'LP1912'),
    @('06:44:15', '07:01', '16_SANTA ANA', 17, 'LP1912'),
    @('05:31:23', '07:04', '23_HERNANDEZ', 93, 'LP1912'),
    @('06:44:15', '07:05', '23_HERNANDEZ', 21, 'LP1912'),
    @('05:31:23', '07:05', '15_ABASTO', 94, 'LP1912'),
    @('05:31:23', '07:07', '225_GOMEZ', 96, 'LP1912'),
    @('05:31:23', '07:11', '215A_EL PATO', 100, 'LP1912'),
    @('05:31:23', '07:15', '11_ETCHEVERRY', 104, 'LP1912'),
    @('06:59:44', '07:16', '16_SANTA ANA', 17, 'LP1912'),
    @('06:44:15', '07:16', '11_ETCHEVERRY', 32, 'LP1912'),
    @('06:44:15', '07:17', '16_SANTA ANA', 33, 'LP1912'),
    @('05:31:23', '07:21', '26_HERNANDEZ', 110, 'LP1912'),
    @('05:31:23', '07:23', '10_OLMOS', 112, 'LP1912'),
    @('06:59:44', '07:24', '10_OLMOS', 25, 'LP1912'),
    @('06:44:15', '07:25', '10_OLMOS', 41, 'LP1912'),
    @('06:01:37', '07:31', '11_ETCHEVERRY', 90, 'LP1912'),
    @('06:01:37', '07:31', '16_SANTA ANA', 90, 'LP1912'),
    @('07:31:43', '07:31', '10_OLMOS', 0, 'LP1912'),
    @('06:44:15', '07:32', '11_ETCHEVERRY', 48, 'LP1912'),
    @('06:01:37', '07:32', '84_COLONIA URQUIZA-ESC 49', 91, 'LP1912'),
    @('06:01:37', '07:36', '27_EL RETIRO', 95, 'LP1912'),
    @('06:44:15', '07:37', '27_EL RETIRO', 53, 'LP1912'),
    @('06:01:37', '07:39', '10_OLMOS', 98, 'LP1912'),
    @('07:31:43', '07:43', '11_ETCHEVERRY', 12, 'LP1912'),
    @('06:59:44', '07:45', '11_ETCHEVERRY', 46, 'LP1912'),
    @('06:01:37', '07:47', '14_ABASTO', 106, 'LP1912'),
    @('06:44:15', '07:48', '14_ABASTO', 64, 'LP1912'),
    @('06:01:37', '07:51', '215D_EL PATO', 110, 'LP1912'),
    @('07:31:43', '07:55', '14_ABASTO', 24, 'LP1912'),
    @('07:31:43', '07:57', '16_SANTA ANA', 26, 'LP1912'),
    @('07:31:43', '07:59', '23_HERNANDEZ', 28, 'LP1912'),
    @('07:31:43', '08:00', '23_HERNANDEZ', 29, 'LP1912'),
    @('07:31:43', '08:03', '11_ETCHEVERRY', 32, 'LP1912'),
    @('06:44:15', '08:04', '23_HERNANDEZ', 80, 'LP1912'),
    @('06:59:44', '08:05', '23_HERNANDEZ', 66, 'LP1912'),
    @('06:44:15', '08:12', '15_ABASTO', 88, 'LP1912'),
    @('06:44:15', '08:21', '26_HERNANDEZ', 97, 'LP1912'),
    @('06:59:44', '08:22', '16_P MOR-SANTA ANA', 83, 'LP1912'),
    @('06:44:15', '08:23', '215B_EL PATO', 99, 'LP1912'),
    @('06:44:15', '08:23', '16_P MOR-SANTA ANA', 99, 'LP1912'),
    @('06:44:15', '08:27', '84_COLONIA URQUIZA-ESC 49', 103, 'LP1912'),
    @('06:44:15', '08:42', '81_EL PELIGRO', 118, 'LP1912'),
    @('07:31:43', '08:43', '14_ABASTO', 72, 'LP1912'),
    @('07:31:43', '08:52', '81_EL PELIGRO', 81, 'LP1912'),
    @('06:59:44', '08:54', '17_ROMERO', 115, 'LP1912'),
    @('07:31:43', '09:01', '215A_EL PATO', 90, 'LP1912'),
    @('07:31:43', '09:10', '16_P MOR-SANTA ANA', 99, 'LP1912'),
    @('07:31:43', '09:16', '27_EL RETIRO', 105, 'LP1912'),
    @('07:31:43', '09:21', '26_HERNANDEZ', 110, 'LP1912'),
    @('07:31:43', '09:22', '17_ROMERO', 111, 'LP1912'),
    @('07:31:43', '09:23', '11_ETCHEVERRY', 112, 'LP1912')
)
for ($i = 0; $i -lt $rows1.Count; $i++) {
    $rowData = $rows1[$i]
    $r = 6 + $i
    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
    $ws.Cells.Item($r, 5).Value = $rowData[4]
}

# --- Sheet 2: LP1912-215 ---
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 1).Value = 'Última actualización: 07:31:43'
$ws.Cells.Item(3, 1).Value = 'Total filas: 10'

$rows2 = @(
    @('03:58:57', '04:45', '215A_EL PATO', 47, 'LP1912'),
    @('04:28:33', '04:46', '215A_EL PATO', 18, 'LP1912'),
    @('03:58:57', '05:34', '215B_EL PATO', 96, 'LP1912'),
    @('05:31:23', '05:35', '215B_EL PATO', 4, 'LP1912'),
    @('04:28:33', '06:11', '215A_EL PATO', 103, 'LP1912'),
    @('04:48:44', '06:46', '215C_EL PATO', 118, 'LP1912'),
    @('05:31:23', '07:11', '215A_EL PATO', 100, 'LP1912'),
    @('06:01:37', '07:51', '215D_EL PATO', 110, 'LP1912'),
    @('06:44:15', '08:23', '215B_EL PATO', 99, 'LP1912'),
    @('07:31:43', '09:01', '215A_EL PATO', 90, 'LP1912')
)
for ($i = 0; $i -lt $rows2.Count; $i++) {
    $rowData = $rows2[$i]
    $r = 6 + $i
    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
    $ws.Cells.Item($r, 5).Value = $rowData[4]
}

# --- Sheet 3: 6203-6173 ---
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 1).Value = 'Última actualización: 07:31:43'
$ws.Cells.Item(3, 1).Value = 'Total filas: 17'

$rows3 = @(
    @('03:58:57', '05:43', '215A_LA PLATA', 105, 'L6173'),
    @('04:28:33', '05:44', '215A_LA PLATA', 76, 'L6173'),
    @('04:48:44', '06:08', '215A_LA PLATA', 80, 'L6173'),
    @('04:28:33', '06:09', '215A_LA PLATA', 101, 'L6173'),
    @('04:48:44', '06:32', '215C_LA PLATA', 104, 'L6203'),
    @('04:58:31', '06:33', '215C_LA PLATA', 95, 'L6203'),
    @('05:31:23', '07:00', '215B_LP-P MOR-1 Y 57', 89, 'L6173'),
    @('06:01:37', '07:35', '215A_LA PLATA', 94, 'L6173'),
    @('06:44:15', '07:36', '215A_LA PLATA', 52, 'L6173'),
    @('06:59:44', '07:39', '215A_LA PLATA', 40, 'L6173'),
    @('07:31:43', '07:40', '215A_LA PLATA', 9, 'L6173'),
    @('06:44:15', '08:07', '215C_LA PLATA', 83, 'L6203'),
    @('07:31:43', '08:08', '215C_LA PLATA', 37, 'L6203'),
    @('06:59:44', '08:09', '215C_LA PLATA', 70, 'L6203'),
    @('06:44:15', '08:31', '215A_LA PLATA', 107, 'L6173'),
    @('06:59:44', '08:35', '215A_LA PLATA', 96, 'L6173'),
    @('07:31:43', '09:08', '215D_LA PLATA', 97, 'L6203')
)
for ($i = 0; $i -lt $rows3.Count; $i++) {
    $rowData = $rows3[$i]
    $r = 6 + $i
    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
    $ws.Cells.Item($r, 5).Value = $rowData[4]
}
